# team_member.xlsx: rename the "电话" (Phone) column header to "分机号"
# (Extension number) on the "member" sheet, and move the active
# selection from D6 to I10 (matching the saved sheet view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: G1 currently reads "电话" -> change to "分机号"
$ws.Range("G1").Value = "分机号"

# Update the saved selection/active cell for the sheet view
$ws.Range("I10").Select()
